# Rename table header cells across the database example slides:
#   "date"   -> "_date"
#   "return" -> "ret"
#   "name"   -> "_name"
#
# The same pair of tables (one with id/ticker/date/return/price,
# one with ticker/name/sector) is repeated on several slides, so walk
# every slide/shape, find any table, and patch any header cell whose
# text exactly matches one of the old values.

$p = $ppt.ActivePresentation

$renames = @{
    "date"   = "_date"
    "return" = "ret"
    "name"   = "_name"
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if (-not $sh.HasTable) {
            continue
        }
        $tbl = $sh.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $cell = $tbl.Cell($r, $c)
                $tr = $cell.Shape.TextFrame.TextRange
                $old = $tr.Text
                if ($renames.ContainsKey($old)) {
                    $tr.Text = $renames[$old]
                }
            }
        }
    }
}
